{"js": "// Fix the typo \"Malo\u00eblle\" -> \"Mano\u00eblle\" in the \"Aanpak & Takenverdeling\"\n// paragraph (planning / task-division text), matching the commit\n// \"planning en takenverdeling aangepast\".\nconst body = context.document.body;\n\n// Narrow the search with surrounding context so we hit the exact\n// occurrence (there are other, already-correct, \"Mano\u00eblle\" spellings\n// elsewhere in the document).\nconst results = body.search(\"waarbij Malo\u00eblle de verantwoordelijkheid\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for the target sentence, found \" + results.items.length);\n}\n\nconst target = results.items[0];\nconst fixedText = target.text.replace(\"Malo\u00eblle\", \"Mano\u00eblle\");\ntarget.insertText(fixedText, \"Replace\");\nawait context.sync();\n", "ps1": "# Fix the typo \"Malo\u00eblle\" -> \"Mano\u00eblle\" in the \"Aanpak & Takenverdeling\"\n# paragraph (planning / task-division text), matching the commit\n# \"planning en takenverdeling aangepast\".\n$d = $word.ActiveDocument\n\n# Search the whole document content; Find.Execute narrows $range down to\n# the matched text in place (there are other, already-correct,\n# \"Mano\u00eblle\" spellings elsewhere in the document, so MatchCase + the\n# exact misspelling keep this to the single intended hit).\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"Malo\u00eblle\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find the target text 'Malo\u00eblle'\"\n}\n\n$range.Text = $range.Text.Replace(\"Malo\u00eblle\", \"Mano\u00eblle\")\n"}
